# Update cryptos list: price (D) and 1h volume/change (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.830.69'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.16%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.348.44'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -1.10%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '544.39'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.69'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -3.00%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.524'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -1.75%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.347.35'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.95%  '
$ws.Range("E10").Value = '  +0.20%  '
$ws.Range("E11").Value = '  +2.09%  '
$ws.Range("E12").Value = '  -0.28%  '
$ws.Range("E13").Value = '  +0.14%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '24.65'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -3.01%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.772.90'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.93%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '61.036.42'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.49%  '
$ws.Range("E17").Value = '  -1.09%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.357.32'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.51%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.63'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.34%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '319.37'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.38%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.11'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.62%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.56'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.85%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.37'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +1.17%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.66'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -8.37%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.38'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +8.74%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.17%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.02'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.77%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '496.97'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -4.51%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.37'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -3.34%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0₃0861'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -6.78%  '
$ws.Range("E32").Value = '  +1.23%  '
$ws.Range("E33").Value = '  -2.09%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.50'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -2.88%  '
$ws.Range("E35").Value = '  +0.08%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.63'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.31%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.376'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.47%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.50'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +3.27%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.23'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -4.17%  '
$ws.Range("E40").Value = '  +5.63%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '143.33'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +5.02%  '
$ws.Range("E42").Value = '  -0.06%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '142.58'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +2.47%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.56'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.86%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.02'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -8.64%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0514'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.11%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '19.11'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -5.49%  '
$ws.Range("E48").Value = '  -1.06%  '
$ws.Range("E49").Value = '  -0.77%  '
$ws.Range("E50").Value = '  -1.46%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '11.39'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.38%  '
